# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 4 (pushing the existing
# records for this market/product down by one row) and populate it
# with the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 4:25 down to 5:26, creating a blank row 4.
$ws.Rows(4).Insert()

# Populate the new row 4 with the latest weekly observation.
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C4").Value = "Arica y Parinacota"
$ws.Range("D4").Value = 45092
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = 100114007
$ws.Range("G4").Value = "Jengibre"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 600
$ws.Range("K4").Value = 13000
$ws.Range("L4").Value = 14000
$ws.Range("M4").Value = 13500
$ws.Range("N4").Value = "$/caja 13 kilos"
$ws.Range("O4").Value = "Perú"
$ws.Range("P4").Value = 1038
$ws.Range("Q4").Value = 13
$ws.Range("R4").Value = "Hortaliza"
